$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing non-availability dates (rows 4-7) ---
$ws.Range("B4").Value = (Get-Date -Year 2018 -Month 4 -Day 18).Date
$ws.Range("B5").Value = (Get-Date -Year 2018 -Month 4 -Day 22).Date

$ws.Range("B6").Value = (Get-Date -Year 2018 -Month 4 -Day 23).Date
# Row 6 gets reassigned to a different VLJ - update the name first, then the
# CSS ID, so the new shared strings land in that order.
$ws.Range("G6").Value = "Morrigan, Jamess"
$ws.Range("F6").Value = "BVAMORRIGA"

$ws.Range("B7").Value = (Get-Date -Year 2018 -Month 4 -Day 25).Date

# --- Fill in the previously-blank row 8 with a new hearing entry ---
$ws.Range("E7").Copy($ws.Range("E8"))
$ws.Range("F7").Copy($ws.Range("F8"))

$ws.Cells.Item(8, 2).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(8, 2).Value = (Get-Date -Year 2018 -Month 4 -Day 26).Date
$ws.Range("C8").Value = "Virtual"
$ws.Range("F8").Value = "BVANULLIGAM"
$ws.Range("H8").Value = 1

# --- Widen the VLJ column to fit the new, longer name ---
$ws.Columns("G").ColumnWidth = 20.5

# --- Move the selection cursor like the author left it ---
$ws.Range("H21").Select() | Out-Null
